$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 686.1818
$ws.Range("I18").Value = 744.7778
$ws.Range("J18").Value = 422.5
$ws.Range("K18").Value = 744.7778
$ws.Range("L18").Value = 422.5
$ws.Range("M18").Value = -460.7778
$ws.Range("N18").Value = -990.5

$ws.Range("H32").Value = 745.82355
$ws.Range("J32").Value = 734.1429000000001
$ws.Range("L32").Value = 734.1429000000001
$ws.Range("N32").Value = -1386.1429

$ws.Range("H51").Value = 2025
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2025
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2025
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -2993

$ws.Range("H64").Value = 2952.9412
$ws.Range("I64").Value = 2800
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 2800
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2552
$ws.Range("N64").Value = -3496

$ws.Range("H67").Value = 2952.9412
$ws.Range("I67").Value = 2800
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 2800
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -1942
$ws.Range("N67").Value = -4716

$ws.Range("H129").Value = 1273.2307
$ws.Range("I129").Value = 328
$ws.Range("J129").Value = 1745.8462
$ws.Range("K129").Value = 984
$ws.Range("L129").Value = 5237.5386
$ws.Range("M129").Value = 4016
$ws.Range("N129").Value = -15237.5386

$ws.Range("H131").Value = 2254.2307
$ws.Range("I131").Value = 470
$ws.Range("J131").Value = 2789.5
$ws.Range("K131").Value = 1410
$ws.Range("L131").Value = 8368.5
$ws.Range("M131").Value = 3630
$ws.Range("N131").Value = -18448.5

$ws.Range("H138").Value = 2965.6711
$ws.Range("I138").Value = 2248.0312
$ws.Range("J138").Value = 3525.7805
$ws.Range("K138").Value = 6744.0936
$ws.Range("L138").Value = 10577.3415
$ws.Range("M138").Value = -1604.0936
$ws.Range("N138").Value = -20857.3415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 995907.8
$ws.Range("I32").Value = 1088328.8
$ws.Range("K32").Value = 1088328.8
$ws.Range("M32").Value = -1088041.8

$ws.Range("H74").Value = 13891673
$ws.Range("I74").Value = 1287.9412
$ws.Range("K74").Value = 1287.9412
$ws.Range("M74").Value = -413.9412

$ws.Range("H77").Value = 13891673
$ws.Range("I77").Value = 1287.9412
$ws.Range("K77").Value = 6439.706
$ws.Range("M77").Value = -2071.706

$ws.Range("H132").Value = 1330223.4
$ws.Range("I132").Value = 3413.0667
$ws.Range("J132").Value = 5923028.5
$ws.Range("K132").Value = 10239.2001
$ws.Range("L132").Value = 17769085.5
$ws.Range("M132").Value = -7709.2001
$ws.Range("N132").Value = -17774145.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6898.16
$ws.Range("I31").Value = 1771.1818
$ws.Range("J31").Value = 8344.23
$ws.Range("K31").Value = 1771.1818
$ws.Range("L31").Value = 8344.23
$ws.Range("M31").Value = -1476.1818
$ws.Range("N31").Value = -8934.23

$ws.Range("H34").Value = 6898.16
$ws.Range("I34").Value = 1771.1818
$ws.Range("J34").Value = 8344.23
$ws.Range("K34").Value = 1771.1818
$ws.Range("L34").Value = 8344.23
$ws.Range("M34").Value = -1569.1818
$ws.Range("N34").Value = -8748.23

$ws.Range("H41").Value = 8650.571
$ws.Range("I41").Value = 5528.5
$ws.Range("J41").Value = 9899.4
$ws.Range("K41").Value = 5528.5
$ws.Range("L41").Value = 9899.4
$ws.Range("M41").Value = -5100.5
$ws.Range("N41").Value = -10755.4

$ws.Range("H50").Value = 12999.2
$ws.Range("J50").Value = 12999.2
$ws.Range("L50").Value = 12999.2
$ws.Range("N50").Value = -14249.2

$ws.Range("H59").Value = 21779.166
$ws.Range("J59").Value = 21779.166
$ws.Range("L59").Value = 21779.166
$ws.Range("N59").Value = -24069.166

$ws.Range("H60").Value = 12341.909
$ws.Range("I60").Value = 5093
$ws.Range("K60").Value = 5093
$ws.Range("M60").Value = -4582

$ws.Range("H68").Value = 22888.334
$ws.Range("J68").Value = 22888.334
$ws.Range("L68").Value = 22888.334
$ws.Range("N68").Value = -24386.334

$ws.Range("H71").Value = 22888.334
$ws.Range("J71").Value = 22888.334
$ws.Range("L71").Value = 68665.00199999999
$ws.Range("N71").Value = -76153.00199999999

$ws.Range("H74").Value = 18467.5
$ws.Range("J74").Value = 18467.5
$ws.Range("L74").Value = 18467.5
$ws.Range("N74").Value = -20215.5

$ws.Range("H77").Value = 18467.5
$ws.Range("J77").Value = 18467.5
$ws.Range("L77").Value = 55402.5
$ws.Range("N77").Value = -64138.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 618.8461
$ws.Range("I5").Value = 423.45456
$ws.Range("J5").Value = 1693.5
$ws.Range("K5").Value = 1270.36368
$ws.Range("L5").Value = 5080.5
$ws.Range("M5").Value = -1158.36368
$ws.Range("N5").Value = -5304.5

$ws.Range("H135").Value = 618.8461
$ws.Range("I135").Value = 423.45456
$ws.Range("J135").Value = 1693.5
$ws.Range("K135").Value = 3811.09104
$ws.Range("L135").Value = 15241.5
$ws.Range("M135").Value = -1276.09104
$ws.Range("N135").Value = -20311.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7598.25
$ws.Range("I70").Value = 7956.7646
$ws.Range("J70").Value = 5566.6665
$ws.Range("K70").Value = 7956.7646
$ws.Range("L70").Value = 5566.6665
$ws.Range("M70").Value = -7686.7646
$ws.Range("N70").Value = -6106.6665

$ws.Range("H73").Value = 7598.25
$ws.Range("I73").Value = 7956.7646
$ws.Range("J73").Value = 5566.6665
$ws.Range("K73").Value = 7956.7646
$ws.Range("L73").Value = 5566.6665
$ws.Range("M73").Value = -7020.7646
$ws.Range("N73").Value = -7438.6665

$ws.Range("H80").Value = 590112.2
$ws.Range("I80").Value = 3001935
$ws.Range("J80").Value = 33537.69
$ws.Range("K80").Value = 3001935
$ws.Range("L80").Value = 33537.69
$ws.Range("M80").Value = -3000937
$ws.Range("N80").Value = -35533.69

$ws.Range("H83").Value = 590112.2
$ws.Range("I83").Value = 3001935
$ws.Range("J83").Value = 33537.69
$ws.Range("K83").Value = 15009675
$ws.Range("L83").Value = 167688.45
$ws.Range("M83").Value = -15004683
$ws.Range("N83").Value = -177672.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 282.32257
$ws.Range("I55").Value = 168.61539
$ws.Range("J55").Value = 364.44446
$ws.Range("K55").Value = 168.61539
$ws.Range("L55").Value = 364.44446
$ws.Range("M55").Value = 4.384610000000009
$ws.Range("N55").Value = -710.4444599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 49032.5
$ws.Range("J135").Value = 49032.5
$ws.Range("L135").Value = 49032.5
$ws.Range("N135").Value = -59172.5
